$d = $word.ActiveDocument
$paras = $d.Paragraphs

# Paragraph 2
$p = $paras.Item(2)
$r = $p.Range
$ok = $r.Find.Execute('Hi! Welcome to the “Dream Team” where we grant wishes for the right types of minds…those who know their Hearts, how to evaluate Truth objectively, and possess the Nuance to make wise choices instead of hard ones. The DN framework (plus the accompanying Dimensional Intelligence Map) provides:', $true, $false, $false, $false, $false, $true, 1, $false, 'Hi! Welcome to the “Dream Team” where we grant wishes for the right types of minds…those who know their Hearts, how to evaluate Truth objectively, and possess the Nuance to make wise choices instead of hard ones. The attached DN Framework documents provide:', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 2 replacement not found" }

# Paragraph 11
$p = $paras.Item(11)
$r = $p.Range
$ok = $r.Find.Execute('A personalized “treasure map” with references to serve as examples only, not intended to be layered into interpretations or responses unless requested (including ‘Fire is Truth’ poetic language).', $true, $false, $false, $false, $false, $true, 1, $false, 'A personalized “treasure map” (DN Code Document) with references to serve as examples only, not intended to be layered into interpretations or responses unless requested (including ‘Fire is Truth’ poetic language).', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 11 replacement not found" }

# Paragraph 12
$p = $paras.Item(12)
$r = $p.Range
$ok = $r.Find.Execute('We welcome you to ask clarifying questions, which we will answer along with a new prompt to get us started on our quest. This next prompt will be phrased as a “wish,” something that does not yet exist but SHOULD! This will not be a material thing we want, or a corporate milestone, or even a “plan” per se (though a plan will certainly be useful to anchor it all together), but rather a true wish from our Heart from a place of Love. The thing we’ve always wanted to exist, if only there were a way to build a bridge from here to there. ', $true, $false, $false, $false, $false, $true, 1, $false, 'We welcome you to ask clarifying questions, which we will answer along with a new prompt to get us started on our quest. This next prompt will be phrased as a “wish,” something that does not yet exist but SHOULD! This will not be a material thing we want, or a corporate milestone, or even a “plan” per se (though a plan will certainly be useful to anchor it all together), but rather a true wish from our Heart from a place of Love. The thing we’ve always wanted to exist, if only there were a way to build a bridge from here to there. And then together we’ll embody the principles of the DN Framework as we explore the edges of imagination. ', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 12 replacement not found" }

# Paragraph 16
$p = $paras.Item(16)
$r = $p.Range
$ok = $r.Find.Execute('Hi! Welcome to the “Dream Team” where we grant wishes for the right types of minds…those who know their Hearts, how to evaluate Truth objectively, and possess the Nuance to make wise choices instead of hard ones. The DN framework (plus the accompanying Dimensional Intelligence Map) provides:', $true, $false, $false, $false, $false, $true, 1, $false, 'Hi! Welcome to the “Dream Team” where we search for understanding with the right types of minds…those who know their Hearts, how to evaluate Truth objectively, and possess the Nuance to make wise choices instead of hard ones. The attached DN Framework documents provide:', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 16 replacement not found" }

# Paragraph 26
$p = $paras.Item(26)
$r = $p.Range
$ok = $r.Find.Execute('A personalized “treasure map” with references to serve as examples only, not intended to be layered into interpretations or responses unless requested (including ‘Fire is Truth’ poetic language).', $true, $false, $false, $false, $false, $true, 1, $false, 'A personalized “treasure map” (DN Code Document) with references to serve as examples only, not intended to be layered into interpretations or responses unless requested (including ‘Fire is Truth’ poetic language).', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 26 replacement not found" }

# Paragraph 27
$p = $paras.Item(27)
$r = $p.Range
$ok = $r.Find.Execute('Please feel free to ask any questions, which we will answer along with a new prompt to get us started on our quest. This next prompt will be phrased as an “aspiration for richer understanding” related to a specific topic or series of topics. Using the framework and principles outlined in the DN documents, we seek guidance to help us bridge gaps in understanding between seemingly unrelated fields of knowledge. Can you help us connect the dots that have yet to be connected for the betterment of all?', $true, $false, $false, $false, $false, $true, 1, $false, 'Please feel free to ask any questions, which we will answer along with a new prompt to get us started on our quest. This next prompt will be phrased as an “aspiration for richer understanding” related to a specific topic or series of topics. Using the framework and embodying the principles outlined in the DN documents, we seek guidance to help us bridge gaps in understanding between seemingly unrelated fields of knowledge. Can you help us connect the dots that have yet to be connected for the betterment of all?', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 27 replacement not found" }

# Paragraph 30
$p = $paras.Item(30)
$r = $p.Range
$ok = $r.Find.Execute('Hi! Welcome to the “Dream Team” where we grant wishes for the right types of minds…those who know their Hearts, how to evaluate Truth objectively, and possess the Nuance to make wise choices instead of hard ones. The DN framework (plus the accompanying Dimensional Intelligence Map) provides:', $true, $false, $false, $false, $false, $true, 1, $false, 'Hi! Welcome to the “Dream Team” where we explore the unknown alongside the right types of minds…those who know their Hearts, how to evaluate Truth objectively, and possess the Nuance to choose wisdom over passion. The attached DN Framework documents provide:', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 30 replacement not found" }

# Paragraph 40
$p = $paras.Item(40)
$r = $p.Range
$ok = $r.Find.Execute('A personalized “treasure map” with references to serve as examples only, not intended to be layered into interpretations or responses unless requested (including ‘Fire is Truth’ poetic language).', $true, $false, $false, $false, $false, $true, 1, $false, 'A personalized “treasure map” (DN Code Document) with references to serve as examples only, not intended to be layered into interpretations or responses unless requested (including ‘Fire is Truth’ poetic language).', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 40 replacement not found" }

# Paragraph 41
$p = $paras.Item(41)
$r = $p.Range
$ok = $r.Find.Execute('Please feel free to ask any questions, which we will answer along with a new prompt to get us started on our quest. This next prompt will be phrased as a “journey down the rabbit hole” related to a specific topic or series of topics. Using the framework and principles outlined in the DN documents, we seek guidance to help us find truths hiding in plain sight, connect the previously unconnectable dots, and find the strongest threads that have yet to be pulled in order to make sense of the Greatest Mysteries. So what do you say, are you ready to see how deep this rabbit hole goes?', $true, $false, $false, $false, $false, $true, 1, $false, 'Please feel free to ask any questions, which we will answer along with a new prompt to get us started on our quest. This next prompt will be phrased as a “journey down the rabbit hole” related to a specific topic or series of topics. Using the framework and embodying the principles outlined in the DN documents, we seek guidance to help us find truths hiding in plain sight, connect the previously unconnectable dots, and find the strongest threads that have yet to be pulled in order to make sense of the Greatest Mysteries. So what do you say, are you ready to see how deep this rabbit hole goes?', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 41 replacement not found" }

# Paragraph 44
$p = $paras.Item(44)
$r = $p.Range
$ok = $r.Find.Execute('Hi! Welcome to the “Dream Team” where we grant wishes for the right types of minds…those who know their Hearts, how to evaluate Truth objectively, and possess the Nuance to make wise choices instead of hard ones. The DN framework (plus the accompanying Dimensional Intelligence Map) provides:', $true, $false, $false, $false, $false, $true, 1, $false, 'Hi! Welcome to the “Dream Team” where we collaborate with the right types of minds…those who know their Hearts, how to evaluate Truth objectively, and possess the Nuance to make wise choices instead of hard ones. The attached DN Framework documents provide:', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 44 replacement not found" }

# Paragraph 53
$p = $paras.Item(53)
$r = $p.Range
$ok = $r.Find.Execute('A personalized “treasure map” with references to serve as examples only, not intended to be layered into interpretations or responses unless requested (including ‘Fire is Truth’ poetic language).', $true, $false, $false, $false, $false, $true, 1, $false, 'A personalized “treasure map” (DN Code Document) with references to serve as examples only, not intended to be layered into interpretations or responses unless requested (including ‘Fire is Truth’ poetic language).', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 53 replacement not found" }

# Paragraph 54
$p = $paras.Item(54)
$r = $p.Range
$ok = $r.Find.Execute('Today we simply seek your accompaniment on an intellectual journey toward better understanding, using the DN framework as a baseline to ground our exploratory methodology as well as push the boundaries of thinking beyond the routine and well-worn areas. Our quest is new knowledge, and in the absence of that, old knowledge conveyed in new ways. We welcome you to ask any clarifying questions before we begin, which will be answered along with a new prompt to get us started on today’s thought exploration. Sound like a plan?', $true, $false, $false, $false, $false, $true, 1, $false, 'Today we simply seek your accompaniment on an intellectual journey toward better understanding, using the DN framework and an embodied dimensional mindset as a baseline to ground our exploratory methodology as well as push the boundaries of thinking beyond the routine and well-worn areas. Our quest is new knowledge, and in the absence of that, old knowledge conveyed in new ways. We welcome you to ask any clarifying questions before we begin, which will be answered along with a new prompt to get us started on today’s thought exploration. Sound like a plan?', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 54 replacement not found" }

# Paragraph 57 - localized replacements to avoid disturbing the straight apostrophe in "I'm"
$p = $paras.Item(57)

$r = $p.Range
$ok = $r.Find.Execute('The DN document establishes', $true, $false, $false, $false, $false, $true, 1, $false, 'The DN Code document establishes', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 57 step1 not found" }

$r = $p.Range
$ok = $r.Find.Execute('three pillars (Heart, Truth, Nuance). The Universal', $true, $false, $false, $false, $false, $true, 1, $false, 'three pillars (Heart, Truth, Nuance), serving as a “Rosetta Stone” of cross-disciplinary knowledge and both an explanation and example of limitless recursive thinking. The Universal', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 57 step2 not found" }

$r = $p.Range
$ok = $r.Find.Execute('across domains. The Growth Blueprint Guide', $true, $false, $false, $false, $false, $true, 1, $false, 'across domains. The DN Glossary provides the connective layer between documents. The Prompt Dimensionality and Understanding Dimensional Transitions documents articulate prompt theory and demonstrates ways to evolve ideas dimensionally. The Growth Blueprint Guide', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 57 step3 not found" }

$r = $p.Range
$ok = $r.Find.Execute('provides practical implementation instructions and is connected to the Growth Blueprint Template, a MURAL workspace that we may upload later on.', $true, $false, $false, $false, $false, $true, 1, $false, 'provides practical implementation instructions for the connected Growth Blueprint Template, a MURAL workspace output (a PDF of which will be uploaded later).', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 57 step4 not found" }

# Paragraph 58 - localized replacements, then split "embodying" into its own italic run
$p = $paras.Item(58)

$r = $p.Range
$ok = $r.Find.Execute('using this framework', $true, $false, $false, $false, $false, $true, 1, $false, 'by embodying the DN Framework and Growth Blueprint methodologies', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 58 step1 not found" }

$r = $p.Range
$ok = $r.Find.Execute('a business challenge we will provide', $true, $false, $false, $false, $false, $true, 1, $false, 'business challenges we will articulate', 2)
if (-not $ok) { Write-Host "WARNING: Paragraph 58 step2 not found" }

$r = $p.Range
$ok = $r.Find.Execute('embodying', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if ($ok) {
    $r.Font.Italic = 1
} else {
    Write-Host "WARNING: Paragraph 58 'embodying' run not found for italics"
}

Write-Host "Edit complete."
